# Leaderboard: correct the existing row 2 battery + add two more lap-time rows,
# per commit "added somr more times".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leaderboard")

# --- Fix existing row 2 (battery was "3S", should be "2S"; new date/time + laptime) ---
$ws.Range("C2").Value = "2S"
$ws.Range("D2").Value = 45984.651388888888
$ws.Range("E2").Value = 19.63

# --- New row 3 ---
$ws.Range("A3").Value = "Charlie"
$ws.Range("B3").Value = "Trophy Course"
$ws.Range("C3").Value = "2S"
$ws.Range("D3").Value = 45984.638194444444
$ws.Range("E3").Value = 18.8

# --- New row 4 ---
$ws.Range("A4").Value = "Charlie"
$ws.Range("B4").Value = "Trophy Course"
$ws.Range("C4").Value = "2S"
$ws.Range("D4").Value = 45984.61041666667
$ws.Range("E4").Value = 19.3

# Date/time columns keep the date+time number format (now the built-in m/d/yy h:mm).
$ws.Range("D2:D4").NumberFormat = "m/d/yy h:mm"

# Column D widened slightly to fit the new date/time values.
$ws.Columns(4).ColumnWidth = 16

# Move the active selection to below the newly-added data, matching the saved view.
[void]$ws.Range("E5").Select()
